# netCrypto.xlsx update - "Add files via upload"
# Bump the USD-amount figure in T2 and move the live selection/scroll
# position from T2 over to Q14 (sheet was re-saved after the author
# scrolled/clicked around near columns L-T before uploading).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the USD Amount value in T2 (256389 -> 257547)
$ws.Range("T2").Value = 257547

# Bring the view roughly into frame (author's sheet was scrolled so column L
# was the left-most visible column) and move the active selection to Q14,
# matching where the cursor was left when the file was saved.
$excel.ActiveWindow.ScrollColumn = 12
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("Q14").Select()
